$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "erledigte Tasks" / progress notes added to the Daily Scrum log (column D),
# plus a couple of absence/task notes in column C, for the entries covering
# the days that were filled in since the last upload.

$ws.Range("D63").Value = "WebService User Fertiggestellt"

$ws.Range("D65").Value = "Webservice Category und Article fertiggestellt"

$ws.Range("D70").Value = "Anbindung Android zu Webservice"

$ws.Range("C75").Value = "Mikula, Unterkofler"
$ws.Range("D75").Value = "Webservice Kommentar fertiggestellt"

$ws.Range("D77").Value = "Fehlerbehebung Android "

$ws.Range("C79").Value = "Bugelnig"
$ws.Range("D79").Value = "User GUI 70% mit Webservice anbindung"

# Restore the view/selection state saved with the workbook (scrolled further
# down the sheet, with D80 as the active cell).
$ws.Range("D80").Select()
